$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The specimen "H 72" (row 2) is dropped from this missing-data sample;
# every row below it shifts up by one (sheet dimension goes from F63 to
# F62). The underlying measurements for every remaining specimen are
# unchanged - only which cells are blanked out (simulated missing data)
# differs between the two random "seed" masks.
$ws.Rows(2).Delete()

# Re-apply the missing-data mask for column B ("A") and column F ("H")
# on the shifted rows so the blanking pattern matches the new seed.
# (Columns A/C/D/E keep the same masking pattern as before the shift.)
$ws.Range("F4").Value = 0.70909
$ws.Range("F8").Value = ""
$ws.Range("F9").Value = 0.71194
$ws.Range("B10").Value = -19.5
$ws.Range("F10").Value = 0.7105
$ws.Range("F11").Value = ""
$ws.Range("B12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("B14").Value = -20.8
$ws.Range("B15").Value = ""
$ws.Range("F20").Value = 0.7106
$ws.Range("F22").Value = ""
$ws.Range("F23").Value = 0.70931
$ws.Range("F25").Value = ""
$ws.Range("B26").Value = -19.5
$ws.Range("B28").Value = ""
$ws.Range("F28").Value = 0.70963
$ws.Range("B30").Value = -19.5
$ws.Range("B32").Value = ""
$ws.Range("F32").Value = ""
$ws.Range("B35").Value = -19.2
$ws.Range("B36").Value = ""
$ws.Range("F36").Value = 0.71087
$ws.Range("B37").Value = -19.8
$ws.Range("B38").Value = ""
$ws.Range("F38").Value = ""
$ws.Range("F41").Value = 0.71115
$ws.Range("F42").Value = 0.71115
$ws.Range("F43").Value = ""
$ws.Range("F44").Value = ""
$ws.Range("B45").Value = -19.7
$ws.Range("B46").Value = ""
$ws.Range("F52").Value = 0.70948
$ws.Range("B53").Value = -20.3
$ws.Range("F54").Value = ""
$ws.Range("B56").Value = ""
